# Remove the empty "ListParagraph" paragraph that immediately follows the
# paragraph containing "צפייה ברשימת ההזמנות" (and precedes the next, also
# empty, paragraph before "הנחות:").
#
# The diff shows a single empty <w:p> (pStyle ListParagraph, indented,
# David font, RTL) being deleted entirely - no text content changes.

$d = $word.ActiveDocument

# Locate the anchor text.
$anchor = $d.Content
$found = $anchor.Find.Execute("צפייה ברשימת ההזמנות", $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The paragraph holding the match text.
    $anchorPara = $anchor.Paragraphs.Item(1)

    # The very next paragraph is the empty one to delete.
    $target = $anchorPara.Range.Next(4, 1)

    # Safety check: only delete it if it is indeed empty (just the pilcrow).
    if (($target.End - $target.Start) -le 1) {
        $target.Delete()
    }
}
